$d = $word.ActiveDocument

# --- Edit 1: Append trailing spaces to the first paragraph's existing text,
#     then append three colored runs forming the red annotation text. ---
$p1 = $d.Paragraphs.Item(1).Range
$p1.End = $p1.End - 1          # exclude the paragraph mark
$p1.InsertAfter("  ")

$run1 = $d.Paragraphs.Item(1).Range
$run1.End = $run1.End - 1
$run1.Collapse(0)               # wdCollapseEnd
$run1.InsertAfter("(This is a change " + [char]0x2013 + " Ve")
$run1.Font.Color = 192          # RGB C00000 -> 0xBBGGRR little-endian long

$run2 = $d.Paragraphs.Item(1).Range
$run2.End = $run2.End - 1
$run2.Collapse(0)
$run2.InsertAfter("rsion for branch alternate")
$run2.Font.Color = 192

$run3 = $d.Paragraphs.Item(1).Range
$run3.End = $run3.End - 1
$run3.Collapse(0)
$run3.InsertAfter(")")
$run3.Font.Color = 192

# --- Edit 2: insert a new, empty, specially-formatted paragraph right after
#     the "It will be treated as a binary file by Git." paragraph. ---
$p2 = $d.Paragraphs.Item(2).Range
$p2.InsertParagraphAfter()

$newParaRange = $d.Paragraphs.Item(3).Range
$newParaRange.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:shd w:val='clear' w:color='auto' w:fill='F9F9F9'/><w:rPr><w:rFonts w:ascii='Calibri' w:eastAsia='Times New Roman' w:hAnsi='Calibri' w:cs='Calibri'/><w:b/><w:bCs/><w:color w:val='202122'/></w:rPr></w:pPr></w:p>")

Write-Host "Edit complete."
